$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.685.72"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "3.772.83"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").Value = "'0.994"
$ws.Range("E4").Value = "  -0.59%  "
$ws.Range("D5").Value = "'598.63"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'162.90"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("D7").Value = "3.770.23"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "'0.156"
$ws.Range("E10").Value = "  -2.92%  "
$ws.Range("D11").Value = "'0.444"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").Value = "'6.62"
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").Value = "4.403.22"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "3.759.03"
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").Value = "67.656.36"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "'18.15"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D20").Value = "'6.98"
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").Value = "'456.30"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("D22").Value = "'9.45"
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("D23").Value = "'0.690"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").Value = "'82.88"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("E25").Value = "  -6.17%  "
$ws.Range("D26").Value = "'11.82"
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'9.83"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "3.915.84"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("E32").Value = "  -2.32%  "
$ws.Range("E33").Value = "  -6.63%  "
$ws.Range("D34").Value = "'28.81"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("E38").Value = "  +3.58%  "
$ws.Range("D39").Value = "'5.75"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'0.977"
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("D41").Value = "'3.16"
$ws.Range("E41").Value = "  -6.44%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D44").Value = "'43.53"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").Value = "'47.15"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("D46").Value = "'151.56"
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "'8.26"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "'384.38"
$ws.Range("E51").Value = "  -2.18%  "
